# Update the crypto price/volume table to reflect the latest scrape.
# Row 39/40 also swap coin identity (LidoDAOToken <-> NEARProtocol).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.053.53"
$ws.Range("E2").Value = "  -1.21%  "

$ws.Range("D3").Value = "2.314.49"
$ws.Range("E3").Value = "  -1.72%  "

$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.65"
$ws.Range("E5").Value = "  -5.43%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "106.48"
$ws.Range("E6").Value = "  +6.51%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.628"
$ws.Range("E7").Value = "  -1.52%  "

$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.611"
$ws.Range("E9").Value = "  -2.93%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.27"
$ws.Range("E10").Value = "  +1.62%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0917"
$ws.Range("E11").Value = "  -0.40%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.31"
$ws.Range("E12").Value = "  -1.68%  "

$ws.Range("E14").Value = "  -2.64%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.58"
$ws.Range("E15").Value = "  -4.45%  "

$ws.Range("D16").Value = "2.653.81"
$ws.Range("E16").Value = "  -2.37%  "

$ws.Range("D17").Value = "2.307.81"
$ws.Range("E17").Value = "  -2.12%  "

$ws.Range("D18").Value = "42.138.21"
$ws.Range("E18").Value = "  -0.98%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.77"
$ws.Range("E19").Value = "  -1.94%  "

$ws.Range("E20").Value = "  -1.54%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "74.71"
$ws.Range("E21").Value = "  -0.80%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.49"
$ws.Range("E22").Value = "  -7.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "258.95"
$ws.Range("E23").Value = "  -3.57%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.29"
$ws.Range("E24").Value = "  -0.86%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.28"
$ws.Range("E25").Value = "  -6.55%  "

$ws.Range("E26").Value = "  +0.47%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.05"
$ws.Range("E27").Value = "  -3.74%  "

$ws.Range("E28").Value = "  +3.30%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.83"
$ws.Range("E29").Value = "  -1.23%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.61"
$ws.Range("E30").Value = "  +0.69%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0895"
$ws.Range("E31").Value = "  -0.75%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "162.58"
$ws.Range("E32").Value = "  -7.40%  "

$ws.Range("E33").Value = "  -5.17%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.88"
$ws.Range("E34").Value = "  -3.45%  "

$ws.Range("E35").Value = "  -2.46%  "

$ws.Range("E36").Value = "  +11.28%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.52"
$ws.Range("E37").Value = "  -1.52%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0354"
$ws.Range("E38").Value = "  -1.26%  "

$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.74"
$ws.Range("E39").Value = "  -7.53%  "

$ws.Range("B40").Value = "NEARProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.64"
$ws.Range("E40").Value = "  -4.60%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "98.50"
$ws.Range("E41").Value = "  +7.91%  "

$ws.Range("E42").Value = "  -3.45%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "70.80"
$ws.Range("E43").Value = "  +0.90%  "

$ws.Range("E44").Value = "  -1.60%  "

$ws.Range("E45").Value = "  -0.17%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.17"
$ws.Range("E46").Value = "  +2.12%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "111.96"
$ws.Range("E47").Value = "  -5.19%  "

$ws.Range("E48").Value = "  -1.62%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.92"
$ws.Range("E49").Value = "  -2.46%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "74.25"
$ws.Range("E50").Value = "  +6.21%  "

$ws.Range("E51").Value = "  +0.06%  "
